$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Set all numeric data cells (D11:H27) to 0
$ws.Range("D11:H14").Value = 0
$ws.Range("D16:H20").Value = 0
$ws.Range("D21:G21").Value = 0
$ws.Range("D22:H22").Value = 0
$ws.Range("D24:H27").Value = 0

# Set the cells that become "-" (dash, text)
$ws.Range("D15:H15").Value = "-"
$ws.Range("H21").Value = "-"
$ws.Range("D23:H23").Value = "-"
